$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "67.834.06"
$ws.Range("E2").Value = "  +2.92%  "
$ws.Range("D3").Value = "2.527.23"
$ws.Range("E3").Value = "  +1.67%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue $ws.Range("D5") "593.14"
$ws.Range("E5").Value = "  +2.63%  "
Set-TextValue $ws.Range("D6") "176.20"
$ws.Range("E6").Value = "  +4.18%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D9").Value = "2.525.41"
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("E10").Value = "  +2.98%  "
$ws.Range("E11").Value = "  +3.18%  "
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("E13").Value = "  -0.21%  "
Set-TextValue $ws.Range("D14") "26.87"
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("D15").Value = "2.982.01"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("E16").Value = "  +3.01%  "
$ws.Range("D17").Value = "67.675.25"
$ws.Range("E17").Value = "  +3.11%  "
$ws.Range("D18").Value = "2.540.77"
$ws.Range("E18").Value = "  +2.29%  "
Set-TextValue $ws.Range("D19") "8.00"
$ws.Range("E19").Value = "  +5.83%  "
Set-TextValue $ws.Range("D20") "11.44"
$ws.Range("E20").Value = "  +3.00%  "
Set-TextValue $ws.Range("D21") "363.10"
$ws.Range("E21").Value = "  +5.95%  "
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("E23").Value = "  +3.02%  "
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D25") "71.22"
$ws.Range("E25").Value = "  +3.57%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D26") "1.00"
$ws.Range("E26").Value = "  -0.08%  "
Set-TextValue $ws.Range("D27") "10.26"
$ws.Range("E27").Value = "  +4.04%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "2.651.56"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").Value = "0.0₃0991"
$ws.Range("E30").Value = "  +3.63%  "
Set-TextValue $ws.Range("D31") "542.71"
$ws.Range("E31").Value = "  +4.38%  "
Set-TextValue $ws.Range("D32") "8.30"
$ws.Range("E32").Value = "  +3.50%  "
$ws.Range("E33").Value = "  +2.87%  "
$ws.Range("E34").Value = "  +3.49%  "
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("E36").Value = "  +0.00%  "
Set-TextValue $ws.Range("D38") "156.52"
$ws.Range("E38").Value = "  -0.67%  "
Set-TextValue $ws.Range("D39") "18.84"
Set-TextValue $ws.Range("D40") "18.67"
$ws.Range("E40").Value = "  +2.01%  "
$ws.Range("E41").Value = "  +1.78%  "
Set-TextValue $ws.Range("D42") "5.19"
$ws.Range("E42").Value = "  +3.37%  "
$ws.Range("E43").Value = "  +2.64%  "
$ws.Range("E44").Value = "  +4.41%  "
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  +2.13%  "
Set-TextValue $ws.Range("D47") "146.54"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("E49").Value = "  +3.11%  "
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("E51").Value = "  +1.39%  "
